# Estadisticos Segundo Parcial 23 Mayo
$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("H2").Value = 8.1

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 5
$ws2.Range("F2").Value = 32
$ws2.Range("G2").Value = 86.48999999999999
$ws2.Range("H2").Value = 8.1

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 17
$ws2.Range("G3").Value = 89.47
$ws2.Range("H3").Value = 7.8

$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 3
$ws2.Range("F4").Value = 25
$ws2.Range("G4").Value = 89.29000000000001
$ws2.Range("H4").Value = 8.699999999999999

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("E2").Value = 5
$ws3.Range("F2").Value = 32
$ws3.Range("G2").Value = 86.48999999999999
$ws3.Range("H2").Value = 8.199999999999999

$ws3.Range("E3").Value = 2
$ws3.Range("F3").Value = 17
$ws3.Range("G3").Value = 89.47
$ws3.Range("H3").Value = 7.9

$ws3.Range("E4").Value = 3
$ws3.Range("F4").Value = 25
$ws3.Range("G4").Value = 89.29000000000001
$ws3.Range("H4").Value = 8.800000000000001

# --- Sheet "Rescatables" ---
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Cells.Item(2, 1).Value = 23330051920225
$ws4.Cells.Item(2, 2).Value = "FLORES"
$ws4.Cells.Item(2, 3).Value = "VAZQUEZ"
$ws4.Cells.Item(2, 4).Value = "MARCO ANTONIO"
$ws4.Cells.Item(2, 5).Value = "DISEÑA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(2, 6).Value = "2AEM"
$ws4.Cells.Item(2, 7).Value = 4

$ws4.Cells.Item(3, 1).Value = 24330051920353
$ws4.Cells.Item(3, 2).Value = "GARCIA"
$ws4.Cells.Item(3, 3).Value = "SANCHEZ"
$ws4.Cells.Item(3, 4).Value = "JOY JARA"
$ws4.Cells.Item(3, 5).Value = "DISEÑA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(3, 6).Value = "2AEM"
$ws4.Cells.Item(3, 7).Value = 4

$ws4.Cells.Item(4, 1).Value = 22330051920006
$ws4.Cells.Item(4, 2).Value = "BAUTISTA"
$ws4.Cells.Item(4, 3).Value = "DE JESUS"
$ws4.Cells.Item(4, 4).Value = "EDER"
$ws4.Cells.Item(4, 5).Value = "REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS"
$ws4.Cells.Item(4, 6).Value = "6AEM"
$ws4.Cells.Item(4, 7).Value = 4

$ws4.Cells.Item(5, 1).Value = 22330051920010
$ws4.Cells.Item(5, 2).Value = "CRUZ"
$ws4.Cells.Item(5, 3).Value = "COYOHUA"
$ws4.Cells.Item(5, 4).Value = "DIEGO"
$ws4.Cells.Item(5, 5).Value = "REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS"
$ws4.Cells.Item(5, 6).Value = "6AEM"
$ws4.Cells.Item(5, 7).Value = 3

$ws4.Cells.Item(6, 1).Value = 22330051920003
$ws4.Cells.Item(6, 2).Value = "APALE"
$ws4.Cells.Item(6, 3).Value = "TEXOCO"
$ws4.Cells.Item(6, 4).Value = "DANIEL"
$ws4.Cells.Item(6, 5).Value = "REALIZA MANTENIMIENTO EN SUBESTACIONES ELÉCTRICAS"
$ws4.Cells.Item(6, 6).Value = "6AEM"
$ws4.Cells.Item(6, 7).Value = 2
